$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) petDefinitions2 table: icon column ([icon], col O) gets a "pet_" prefix ---
# Rows 5-14 hold the baby-dragon definitions; O5:O14 currently store the bare sku
# (e.g. "baby_classic"); rename them to the "pet_" prefixed icon id.
$iconUpdates = @{
    5  = "pet_baby_classic"
    6  = "pet_baby_crocodile"
    7  = "pet_baby_titan"
    8  = "pet_baby_jawfrey"
    9  = "pet_baby_dark"
    10 = "pet_baby_dino"
    11 = "pet_baby_alien"
    12 = "pet_baby_devil"
    13 = "pet_baby_tony"
    14 = "pet_baby_hedgehog"
}
foreach ($row in $iconUpdates.Keys) {
    $ws.Range("O" + $row).Value = $iconUpdates[$row]
}

# --- 2) petCategoryDefinitions3 table (B19:G20) grows two columns: [tidName] & [tidDescription] ---
$lo = $ws.ListObjects.Item("petCategoryDefinitions3")
$lo.ListColumns.Add() | Out-Null
$lo.ListColumns.Add() | Out-Null

# Header row (row 19)
$ws.Range("H19").Value = "[tidName]"
$ws.Range("I19").Value = "[tidDescription]"

# Data row (row 20)
$ws.Range("H20").Value = "TID_SHARED_EXTRA_GEMS_NAME"
$ws.Range("I20").Value = "TID_SHARED_EXTRA_GEMS_DESC"

# Carry over the look of the preceding column ([firstSucceed], G) onto the two
# new columns, for both the header and the data row.
$ws.Range("G19").Copy() | Out-Null
$ws.Range("H19:I19").PasteSpecial(-4122) | Out-Null
$ws.Range("G20").Copy() | Out-Null
$ws.Range("H20:I20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 3) Leave the view the way the author left it: scrolled to row 4, J20 selected ---
$ws.Range("J20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
